$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (data point removed)
$ws.Rows.Item(6).Delete()

# Adjust column widths (stored XML width 7 -> 8, one col 8 -> 9)
$ws.Columns.Item(2).ColumnWidth = 7.166666666666667
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(17).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(29).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666

# Update data rows 2-5 with new sensor readings
$ws.Cells.Item(2, 1).Value = 45082.50694444445
$ws.Cells.Item(2, 2).Value = 7.205
$ws.Cells.Item(2, 3).Value = 5.605
$ws.Cells.Item(2, 4).Value = 1.192
$ws.Cells.Item(2, 5).Value = 15.61
$ws.Cells.Item(2, 6).Value = 12.374
$ws.Cells.Item(2, 7).Value = 4.794
$ws.Cells.Item(2, 8).Value = 14.798
$ws.Cells.Item(2, 9).Value = 8.957000000000001
$ws.Cells.Item(2, 10).Value = 4.429
$ws.Cells.Item(2, 11).Value = 5.631
$ws.Cells.Item(2, 12).Value = 6.249
$ws.Cells.Item(2, 13).Value = 7.306
$ws.Cells.Item(2, 14).Value = 2.788
$ws.Cells.Item(2, 15).Value = 6.015
$ws.Cells.Item(2, 16).Value = 7.738
$ws.Cells.Item(2, 17).Value = 5.138
$ws.Cells.Item(2, 18).Value = 0.492
$ws.Cells.Item(2, 19).Value = 0.931
$ws.Cells.Item(2, 20).Value = 84.14
$ws.Cells.Item(2, 21).Value = 16.424
$ws.Cells.Item(2, 22).Value = 5.552
$ws.Cells.Item(2, 23).Value = 10.174
$ws.Cells.Item(2, 24).Value = 6.283
$ws.Cells.Item(2, 25).Value = 0.894
$ws.Cells.Item(2, 26).Value = 9.741
$ws.Cells.Item(2, 27).Value = 4.361
$ws.Cells.Item(2, 28).Value = 4.885
$ws.Cells.Item(2, 29).Value = 6.06
$ws.Cells.Item(2, 30).Value = 8.134
$ws.Cells.Item(2, 31).Value = 1.522
$ws.Cells.Item(2, 32).Value = 13.18
$ws.Cells.Item(2, 33).Value = 3.625
$ws.Cells.Item(2, 34).Value = 6.387
$ws.Cells.Item(3, 1).Value = 45082.51388888889
$ws.Cells.Item(3, 2).Value = 14.614
$ws.Cells.Item(3, 3).Value = 11.031
$ws.Cells.Item(3, 4).Value = 0.888
$ws.Cells.Item(3, 5).Value = 32.028
$ws.Cells.Item(3, 6).Value = 25.96
$ws.Cells.Item(3, 7).Value = 10.959
$ws.Cells.Item(3, 8).Value = 41.053
$ws.Cells.Item(3, 9).Value = 17.804
$ws.Cells.Item(3, 10).Value = 8.425000000000001
$ws.Cells.Item(3, 11).Value = 11.631
$ws.Cells.Item(3, 12).Value = 12.796
$ws.Cells.Item(3, 13).Value = 13.925
$ws.Cells.Item(3, 14).Value = 4.267
$ws.Cells.Item(3, 15).Value = 11.655
$ws.Cells.Item(3, 16).Value = 16.181
$ws.Cells.Item(3, 17).Value = 9.856999999999999
$ws.Cells.Item(3, 18).Value = 0.272
$ws.Cells.Item(3, 19).Value = 0.788
$ws.Cells.Item(3, 20).Value = 170.062
$ws.Cells.Item(3, 21).Value = 32.52
$ws.Cells.Item(3, 22).Value = 10.758
$ws.Cells.Item(3, 23).Value = 21.503
$ws.Cells.Item(3, 24).Value = 11.804
$ws.Cells.Item(3, 25).Value = 1.558
$ws.Cells.Item(3, 26).Value = 21.894
$ws.Cells.Item(3, 27).Value = 9.202999999999999
$ws.Cells.Item(3, 28).Value = 8.739000000000001
$ws.Cells.Item(3, 29).Value = 10.362
$ws.Cells.Item(3, 30).Value = 14.152
$ws.Cells.Item(3, 31).Value = 0.707
$ws.Cells.Item(3, 32).Value = 37.737
$ws.Cells.Item(3, 33).Value = 6.411
$ws.Cells.Item(3, 34).Value = 13.147
$ws.Cells.Item(4, 1).Value = 45082.52083333334
$ws.Cells.Item(4, 2).Value = 15.177
$ws.Cells.Item(4, 3).Value = 11.418
$ws.Cells.Item(4, 4).Value = 0.76
$ws.Cells.Item(4, 5).Value = 33.251
$ws.Cells.Item(4, 6).Value = 27.07
$ws.Cells.Item(4, 7).Value = 11.55
$ws.Cells.Item(4, 8).Value = 46.445
$ws.Cells.Item(4, 9).Value = 18.445
$ws.Cells.Item(4, 10).Value = 8.644
$ws.Cells.Item(4, 11).Value = 12.107
$ws.Cells.Item(4, 12).Value = 13.3
$ws.Cells.Item(4, 13).Value = 14.336
$ws.Cells.Item(4, 14).Value = 4.24
$ws.Cells.Item(4, 15).Value = 12.031
$ws.Cells.Item(4, 16).Value = 16.865
$ws.Cells.Item(4, 17).Value = 10.142
$ws.Cells.Item(4, 18).Value = 0.204
$ws.Cells.Item(4, 19).Value = 0.6830000000000001
$ws.Cells.Item(4, 20).Value = 175.798
$ws.Cells.Item(4, 21).Value = 33.652
$ws.Cells.Item(4, 22).Value = 11.105
$ws.Cells.Item(4, 23).Value = 22.423
$ws.Cells.Item(4, 24).Value = 12.127
$ws.Cells.Item(4, 25).Value = 1.591
$ws.Cells.Item(4, 26).Value = 23.528
$ws.Cells.Item(4, 27).Value = 9.603
$ws.Cells.Item(4, 28).Value = 8.913
$ws.Cells.Item(4, 29).Value = 10.513
$ws.Cells.Item(4, 30).Value = 14.41
$ws.Cells.Item(4, 31).Value = 0.461
$ws.Cells.Item(4, 32).Value = 42.583
$ws.Cells.Item(4, 33).Value = 6.529
$ws.Cells.Item(4, 34).Value = 13.676
$ws.Cells.Item(5, 1).Value = 45082.52777777778
$ws.Cells.Item(5, 2).Value = 4.65
$ws.Cells.Item(5, 3).Value = 3.5
$ws.Cells.Item(5, 4).Value = 0.34
$ws.Cells.Item(5, 5).Value = 10.33
$ws.Cells.Item(5, 6).Value = 8.19
$ws.Cells.Item(5, 7).Value = 3.35
$ws.Cells.Item(5, 8).Value = 19.03
$ws.Cells.Item(5, 9).Value = 5.68
$ws.Cells.Item(5, 10).Value = 2.9
$ws.Cells.Item(5, 11).Value = 3.58
$ws.Cells.Item(5, 12).Value = 4.11
$ws.Cells.Item(5, 13).Value = 4.55
$ws.Cells.Item(5, 14).Value = 1.51
$ws.Cells.Item(5, 15).Value = 3.76
$ws.Cells.Item(5, 16).Value = 5.18
$ws.Cells.Item(5, 17).Value = 3.27
$ws.Cells.Item(5, 18).Value = 0.09
$ws.Cells.Item(5, 19).Value = 0.31
$ws.Cells.Item(5, 20).Value = 49.94
$ws.Cells.Item(5, 21).Value = 10.76
$ws.Cells.Item(5, 22).Value = 3.47
$ws.Cells.Item(5, 23).Value = 7.04
$ws.Cells.Item(5, 24).Value = 3.91
$ws.Cells.Item(5, 25).Value = 0.5
$ws.Cells.Item(5, 26).Value = 9.460000000000001
$ws.Cells.Item(5, 27).Value = 2.91
$ws.Cells.Item(5, 28).Value = 2.94
$ws.Cells.Item(5, 29).Value = 3.46
$ws.Cells.Item(5, 30).Value = 4.62
$ws.Cells.Item(5, 31).Value = 0.34
$ws.Cells.Item(5, 32).Value = 17.9
$ws.Cells.Item(5, 33).Value = 2.13
$ws.Cells.Item(5, 34).Value = 4.18
